# Regenerate merged AHB files
# - Rename the "_old" / "_new" column headers to "_FV2404" / "_FV2410"
# - Turn the header + data range into an Excel Table ("Table1")
# - Freeze the header row (row 1) in the sheet view

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# 1) Rename the header row (row 1) shared strings in place.
$oldSuffixHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)
$newSuffixHeaders = @(
    "Segmentname_FV2404",
    "Segmentgruppe_FV2404",
    "Segment_FV2404",
    "Datenelement_FV2404",
    "Segment ID_FV2404",
    "Code_FV2404",
    "Qualifier_FV2404",
    "Beschreibung_FV2404",
    "Bedingungsausdruck_FV2404",
    "Bedingung_FV2404"
)

for ($i = 0; $i -lt $oldSuffixHeaders.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $newSuffixHeaders[$i]
}

$oldNewHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)
$newNewHeaders = @(
    "Segmentname_FV2410",
    "Segmentgruppe_FV2410",
    "Segment_FV2410",
    "Datenelement_FV2410",
    "Segment ID_FV2410",
    "Code_FV2410",
    "Qualifier_FV2410",
    "Beschreibung_FV2410",
    "Bedingungsausdruck_FV2410",
    "Bedingung_FV2410"
)

# Columns L..U are indices 12..21 (column K = 11 holds "diff" and is unchanged)
for ($i = 0; $i -lt $oldNewHeaders.Length; $i++) {
    $ws.Cells.Item(1, 11 + 1 + $i).Value = $newNewHeaders[$i]
}

# 2) Turn A1:U60 into a native Excel table named "Table1" with an AutoFilter.
$tableRange = $ws.Range("A1:U60")
$lo = $ws.ListObjects.Add(1, $tableRange, $null, 1)
$lo.Name = "Table1"

# 3) Freeze panes at row 1 (so the header row stays visible while scrolling).
[void]$ws.Range("A2").Select()
[void]($excel.ActiveWindow.FreezePanes = $true)
